# Auto-generated edit script applying the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '60.725.87'
$ws.Range("E2").Value = '  -1.91%  '

# Row 3
$ws.Range("D3").Value = '2.906.92'
$ws.Range("E3").Value = '  -3.05%  '

# Row 4
$ws.Range("E4").Value = '  +0.09%  '

# Row 5
$ws.Range("D5").Value = "'584.62"
$ws.Range("E5").Value = '  -1.27%  '

# Row 6
$ws.Range("D6").Value = "'147.39"
$ws.Range("E6").Value = '  +0.59%  '

# Row 7
$ws.Range("E7").Value = '  -0.10%  '

# Row 8
$ws.Range("B8").Value = 'XRP'
$ws.Range("C8").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D8").Value = "'0.502"
$ws.Range("E8").Value = '  -2.70%  '

# Row 9
$ws.Range("B9").Value = 'LidoStakedEther'
$ws.Range("C9").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D9").Value = '2.900.72'
$ws.Range("E9").Value = '  -3.13%  '

# Row 10
$ws.Range("D10").Value = "'6.70"
$ws.Range("E10").Value = '  +7.71%  '

# Row 11
$ws.Range("D11").Value = "'0.144"
$ws.Range("E11").Value = '  -2.95%  '

# Row 12
$ws.Range("D12").Value = "'0.446"
$ws.Range("E12").Value = '  -2.50%  '

# Row 13
$ws.Range("D13").Value = "'0.0000224"
$ws.Range("E13").Value = '  -2.85%  '

# Row 14
$ws.Range("D14").Value = "'34.34"
$ws.Range("E14").Value = '  +0.32%  '

# Row 15
$ws.Range("D15").Value = "'0.127"

# Row 16
$ws.Range("D16").Value = '3.389.98'
$ws.Range("E16").Value = '  -3.03%  '

# Row 17
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '60.823.39'
$ws.Range("E17").Value = '  -1.77%  '

# Row 18
$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D18").Value = "'6.82"
$ws.Range("E18").Value = '  -2.39%  '

# Row 19
$ws.Range("D19").Value = '2.909.56'
$ws.Range("E19").Value = '  -2.93%  '

# Row 20
$ws.Range("D20").Value = "'424.45"
$ws.Range("E20").Value = '  -4.98%  '

# Row 21
$ws.Range("D21").Value = "'13.61"
$ws.Range("E21").Value = '  -3.60%  '

# Row 22
$ws.Range("D22").Value = "'0.670"
$ws.Range("E22").Value = '  -2.08%  '

# Row 23
$ws.Range("D23").Value = "'7.14"
$ws.Range("E23").Value = '  -2.97%  '

# Row 24
$ws.Range("D24").Value = "'80.88"
$ws.Range("E24").Value = '  -1.32%  '

# Row 25
$ws.Range("D25").Value = "'11.01"
$ws.Range("E25").Value = '  +0.63%  '

# Row 26
$ws.Range("D26").Value = "'2.18"
$ws.Range("E26").Value = '  -1.73%  '

# Row 27
$ws.Range("D27").Value = "'11.80"
$ws.Range("E27").Value = '  -2.45%  '

# Row 28
$ws.Range("E28").Value = '  -0.11%  '

# Row 29
$ws.Range("B29").Value = 'FirstDigitalUSD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = '  -0.03%  '

# Row 30
$ws.Range("B30").Value = 'NEARProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D30").Value = "'7.26"
$ws.Range("E30").Value = '  +0.49%  '

# Row 31
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").Value = "'2.18"
$ws.Range("E31").Value = '  +4.79%  '

# Row 32
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = "'2.62"
$ws.Range("E32").Value = '  -2.86%  '

# Row 33
$ws.Range("D33").Value = "'26.65"
$ws.Range("E33").Value = '  -2.64%  '

# Row 34
$ws.Range("D34").Value = "'0.106"
$ws.Range("E34").Value = '  -3.67%  '

# Row 35
$ws.Range("D35").Value = '0.0₃0838'
$ws.Range("E35").Value = '  -0.66%  '

# Row 36
$ws.Range("E36").Value = '  -1.20%  '

# Row 37
$ws.Range("D37").Value = "'5.67"
$ws.Range("E37").Value = '  -2.33%  '

# Row 38
$ws.Range("B38").Value = 'OKB'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D38").Value = "'49.61"
$ws.Range("E38").Value = '  -0.91%  '

# Row 39
$ws.Range("B39").Value = 'dogwifhat'
$ws.Range("C39").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D39").Value = "'2.98"
$ws.Range("E39").Value = '  +2.57%  '

# Row 40
$ws.Range("D40").Value = "'2.03"
$ws.Range("E40").Value = '  +0.19%  '

# Row 41
$ws.Range("D41").Value = "'0.124"
$ws.Range("E41").Value = '  +1.12%  '

# Row 42
$ws.Range("D42").Value = "'8.73"
$ws.Range("E42").Value = '  -2.50%  '

# Row 43
$ws.Range("D43").Value = "'0.289"
$ws.Range("E43").Value = '  +3.53%  '

# Row 44
$ws.Range("D44").Value = "'41.49"
$ws.Range("E44").Value = '  +1.69%  '

# Row 45
$ws.Range("D45").Value = "'376.94"
$ws.Range("E45").Value = '  -4.97%  '

# Row 46
$ws.Range("D46").Value = "'0.0346"
$ws.Range("E46").Value = '  -1.16%  '

# Row 47
$ws.Range("D47").Value = '2.651.81'
$ws.Range("E47").Value = '  -2.01%  '

# Row 48
$ws.Range("D48").Value = "'132.91"
$ws.Range("E48").Value = '  +0.44%  '

# Row 49
$ws.Range("B49").Value = 'USDe'
$ws.Range("C49").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D49").Value = "'1.00"
$ws.Range("E49").Value = '  +0.02%  '

# Row 50
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").Value = "'25.36"
$ws.Range("E50").Value = '  +7.22%  '

# Row 51
$ws.Range("E51").Value = '  -0.70%  '
